# Ajout du programme et CTE corrigé
# Corrige la coquille dans la description du Playtest sur prototypes papier
# (feuille "tblTypeTest", cellule D10) et remet la sélection active sur C3.

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("tblTypeTest")

# Corrige le texte (fautes d'accord / d'orthographe) dans la cellule D10
$wsTest.Range("D10").Value = "Construire un prototype papier du jeu et le tester avec plusieurs joueurs et un modérateur qui s'assure que les joueurs ne sont pas perdus. Analyser les retours des joueurs et leurs actions en jeu."

# Réactive la feuille "tblTypeTest" et replace la sélection sur C3
$wsTest.Activate()
$wsTest.Range("C3").Select()
